$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# --- Azure parallelism & metrics table cell updates ---
$ws.Range("BG3").Value = 1
$ws.Range("BG4").Value = 1
$ws.Range("BG5").Value = 1
$ws.Range("AP6").Value = 0
$ws.Range("AQ6").Value = 1
$ws.Range("AR6").Value = 0.025
$ws.Range("AS6").Value = 0.03
$ws.Range("AT6").Value = 0
$ws.Range("AU6").Value = 0
$ws.Range("AV6").Value = 0
$ws.Range("AW6").Value = 0
$ws.Range("AX6").Value = 0
$ws.Range("AY6").Value = 0
$ws.Range("AZ6").Value = 0
$ws.Range("BA6").Value = 0
$ws.Range("BG6").Value = 1000
$ws.Range("AP7").Value = 1.5
$ws.Range("AQ7").Value = 1.6
$ws.Range("AR7").Value = 0.027000000000000003
$ws.Range("AT7").Value = 0.008
$ws.Range("AU7").Value = 6
$ws.Range("AV7").Value = 6
$ws.Range("AW7").Value = 0.0015
$ws.Range("AX7").Value = 0.008
$ws.Range("AY7").Value = 6
$ws.Range("AZ7").Value = 6
$ws.Range("BA7").Value = 0.0015
$ws.Range("T8").Value = 2
$ws.Range("AM8").Value = 0.8
$ws.Range("AN8").Value = 0.75
$ws.Range("AP8").Value = 0.75
$ws.Range("AQ8").Value = 1.5
$ws.Range("AR8").Value = 0.025
$ws.Range("AS8").Value = 0.013999999999999999
$ws.Range("AT8").Value = 0.004
$ws.Range("AU8").Value = 6
$ws.Range("AV8").Value = 6
$ws.Range("AW8").Value = 0.001
$ws.Range("AX8").Value = 0.004
$ws.Range("AY8").Value = 6
$ws.Range("AZ8").Value = 6
$ws.Range("BA8").Value = 0.001
$ws.Range("BB8").Value = 0.001
$ws.Range("P9").Value = 32
$ws.Range("AN9").Value = 0.75
$ws.Range("AP9").Value = 1.5
$ws.Range("AQ9").Value = 1.6
$ws.Range("AR9").Value = 0.027000000000000003
$ws.Range("AS9").Value = 0.008
$ws.Range("AT9").Value = 0.008
$ws.Range("AU9").Value = 6
$ws.Range("AV9").Value = 6
$ws.Range("AW9").Value = 0.0015
$ws.Range("AX9").Value = 0.008
$ws.Range("AY9").Value = 6
$ws.Range("AZ9").Value = 6
$ws.Range("BA9").Value = 0.0015
$ws.Range("BB9").Value = 0.0005
$ws.Range("BG9").Value = 1000
$ws.Range("K10").Value = 0.05
$ws.Range("L10").Value = 12
$ws.Range("M10").Value = 12
$ws.Range("N10").Value = 0.013999999999999999
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 12
$ws.Range("Q10").Value = 12
$ws.Range("R10").Value = 0.012
$ws.Range("S10").Value = 0.06
$ws.Range("T10").Value = 0
$ws.Range("V10").Value = 20
$ws.Range("AD10").Value = 0.1275
$ws.Range("AG10").Value = 0.87
$ws.Range("AJ10").Value = 0.029900000000000003
$ws.Range("AK10").Value = 24
$ws.Range("AL10").Value = 240
$ws.Range("AO10").Value = 0.6
$ws.Range("AP10").Value = 2
$ws.Range("AU10").Value = 2
$ws.Range("AV10").Value = 2
$ws.Range("AY10").Value = 2
$ws.Range("AZ10").Value = 2
$ws.Range("BG11").Value = 1
$ws.Range("S12").Value = 0.05
$ws.Range("V12").Value = 25
$ws.Range("AG12").Value = 0.75
$ws.Range("AO12").Value = 0.6
$ws.Range("BG12").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 12
$ws.Range("Q13").Value = 12
$ws.Range("R13").Value = 0.012
$ws.Range("S13").Value = 0.084
$ws.Range("V13").Value = 17
$ws.Range("AD13").Value = 0.15
$ws.Range("AG13").Value = 0.98
$ws.Range("AJ13").Value = 0.0149
$ws.Range("AK13").Value = 12
$ws.Range("AL13").Value = 240
$ws.Range("AO13").Value = 0.45
$ws.Range("AP13").Value = 0.75
$ws.Range("BB13").Value = 0.0001
$ws.Range("BC13").Value = 0
$ws.Range("BD13").Value = 0
$ws.Range("BG13").Value = 1000

# --- Restore selection/pane state to match the authored view ---
$ws.Activate()
[void]$ws.Range("AR6:AR9").Select()
